# Auto-generated edit script: update market-price derived columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets, per scheduled price refresh.
$wb = $excel.ActiveWorkbook

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 36248.832
$ws.Range("I98").Value = 37785.957
$ws.Range("K98").Value = 37785.957
$ws.Range("M98").Value = -36287.957

# ALC row 100
$ws.Range("H100").Value = 1580.4348
$ws.Range("I100").Value = 596.9375
$ws.Range("K100").Value = 596.9375
$ws.Range("M100").Value = -55.9375

# ALC row 108
$ws.Range("H108").Value = 47969.715
$ws.Range("I108").Value = 46800
$ws.Range("J108").Value = 48437.6
$ws.Range("K108").Value = 46800
$ws.Range("L108").Value = 48437.6
$ws.Range("M108").Value = -42960
$ws.Range("N108").Value = -56117.6

# ALC row 113
$ws.Range("H113").Value = 5333.4707
$ws.Range("I113").Value = 4422.5
$ws.Range("J113").Value = 6143.222
$ws.Range("K113").Value = 4422.5
$ws.Range("L113").Value = 6143.222
$ws.Range("M113").Value = -1168.5
$ws.Range("N113").Value = -12651.222

# ALC row 122
$ws.Range("H122").Value = 36248.832
$ws.Range("I122").Value = 37785.957
$ws.Range("K122").Value = 113357.871
$ws.Range("M122").Value = -110907.871

# ALC row 129
$ws.Range("H129").Value = 2809.7856
$ws.Range("J129").Value = 5171.3335
$ws.Range("L129").Value = 15514.0005
$ws.Range("N129").Value = -25514.0005

# ALC row 135
$ws.Range("H135").Value = 13159456
$ws.Range("I135").Value = 1272.3334
$ws.Range("J135").Value = 35716340
$ws.Range("K135").Value = 11451.0006
$ws.Range("L135").Value = 321447060
$ws.Range("M135").Value = -8916.000599999999
$ws.Range("N135").Value = -321452130

# ALC row 138
$ws.Range("H138").Value = 7251566
$ws.Range("I138").Value = 2513.2222
$ws.Range("J138").Value = 9014849
$ws.Range("K138").Value = 7539.6666
$ws.Range("L138").Value = 27044547
$ws.Range("M138").Value = -2399.6666
$ws.Range("N138").Value = -27054827

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 34486540
$ws.Range("I61").Value = 50002460
$ws.Range("K61").Value = 50002460
$ws.Range("M61").Value = -50002248

# ARM row 74
$ws.Range("H74").Value = 62574184
$ws.Range("I74").Value = 125142250
$ws.Range("K74").Value = 125142250
$ws.Range("M74").Value = -125141376

# ARM row 77
$ws.Range("H77").Value = 62574184
$ws.Range("I77").Value = 125142250
$ws.Range("K77").Value = 625711250
$ws.Range("M77").Value = -625706882

# ARM row 110
$ws.Range("H110").Value = 18944.682
$ws.Range("I110").Value = 22917.176
$ws.Range("J110").Value = 5438.2
$ws.Range("K110").Value = 22917.176
$ws.Range("L110").Value = 5438.2
$ws.Range("M110").Value = -20872.176
$ws.Range("N110").Value = -9528.200000000001

# ARM row 136
$ws.Range("H136").Value = 34486540
$ws.Range("I136").Value = 50002460
$ws.Range("K136").Value = 150007380
$ws.Range("M136").Value = -150004830

# BSM row 64
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1871.3334
$ws.Range("J64").Value = 2215.8572
$ws.Range("L64").Value = 2215.8572
$ws.Range("N64").Value = -2665.8572

# BSM row 67
$ws.Range("H67").Value = 1871.3334
$ws.Range("J67").Value = 2215.8572
$ws.Range("L67").Value = 2215.8572
$ws.Range("N67").Value = -3775.8572

# BSM row 99
$ws.Range("H99").Value = 5145.4287
$ws.Range("I99").Value = 3802.3333
$ws.Range("J99").Value = 6152.75
$ws.Range("K99").Value = 3802.3333
$ws.Range("L99").Value = 6152.75
$ws.Range("M99").Value = -2304.3333
$ws.Range("N99").Value = -9148.75

# BSM row 107
$ws.Range("H107").Value = 5341.478
$ws.Range("I107").Value = 3985.5386
$ws.Range("J107").Value = 7104.2
$ws.Range("K107").Value = 3985.5386
$ws.Range("L107").Value = 7104.2
$ws.Range("M107").Value = -2065.5386
$ws.Range("N107").Value = -10944.2

# BSM row 134
$ws.Range("H134").Value = 4685.684
$ws.Range("I134").Value = 4531.0586
$ws.Range("K134").Value = 13593.1758
$ws.Range("M134").Value = -11058.1758

# CRP row 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 688.8182
$ws.Range("I10").Value = 453.55554
$ws.Range("J10").Value = 1747.5
$ws.Range("K10").Value = 453.55554
$ws.Range("L10").Value = 1747.5
$ws.Range("M10").Value = -314.55554
$ws.Range("N10").Value = -2025.5

# CRP row 31
$ws.Range("H31").Value = 38467504
$ws.Range("I31").Value = 5584.467
$ws.Range("J31").Value = 90915576
$ws.Range("K31").Value = 5584.467
$ws.Range("L31").Value = 90915576
$ws.Range("M31").Value = -5289.467
$ws.Range("N31").Value = -90916166

# CRP row 34
$ws.Range("H34").Value = 38467504
$ws.Range("I34").Value = 5584.467
$ws.Range("J34").Value = 90915576
$ws.Range("K34").Value = 5584.467
$ws.Range("L34").Value = 90915576
$ws.Range("M34").Value = -5382.467
$ws.Range("N34").Value = -90915980

# CRP row 103
$ws.Range("H103").Value = 2692.2856
$ws.Range("I103").Value = 2692.2856
$ws.Range("K103").Value = 2692.2856
$ws.Range("M103").Value = -1520.2856

# CRP row 107
$ws.Range("H107").Value = 1219.3889
$ws.Range("I107").Value = 711.0714
$ws.Range("J107").Value = 2998.5
$ws.Range("K107").Value = 711.0714
$ws.Range("L107").Value = 2998.5
$ws.Range("M107").Value = 1208.9286
$ws.Range("N107").Value = -6838.5

# CRP row 122
$ws.Range("H122").Value = 1835.3684
$ws.Range("I122").Value = 1853.6666
$ws.Range("J122").Value = 1818.9
$ws.Range("K122").Value = 5560.9998
$ws.Range("L122").Value = 5456.700000000001
$ws.Range("M122").Value = -3110.9998
$ws.Range("N122").Value = -10356.7

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 36501624
$ws.Range("I4").Value = 55471756
$ws.Range("K4").Value = 166415268
$ws.Range("M4").Value = -166415156

# CUL row 26
$ws.Range("H26").Value = 286.27274
$ws.Range("J26").Value = 557.8
$ws.Range("L26").Value = 1673.4
$ws.Range("N26").Value = -2249.4

# GSM row 58
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 22000
$ws.Range("I58").Value = 22000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 22000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -21723
$ws.Range("N58").ClearContents()

# GSM row 107
$ws.Range("H107").Value = 636.2632
$ws.Range("I107").Value = 690.1667
$ws.Range("J107").Value = 543.8570999999999
$ws.Range("K107").Value = 690.1667
$ws.Range("L107").Value = 543.8570999999999
$ws.Range("M107").Value = 1229.8333
$ws.Range("N107").Value = -4383.8571

# GSM row 114
$ws.Range("H114").Value = 72361
$ws.Range("J114").Value = 72361
$ws.Range("L114").Value = 72361
$ws.Range("N114").Value = -81039

# GSM row 118
$ws.Range("H118").Value = 18820.25
$ws.Range("I118").Value = 20381
$ws.Range("J118").Value = 18300
$ws.Range("K118").Value = 20381
$ws.Range("L118").Value = 18300
$ws.Range("M118").Value = -18724
$ws.Range("N118").Value = -21614

# GSM row 126
$ws.Range("H126").Value = 16673409
$ws.Range("I126").Value = 14295652
$ws.Range("J126").Value = 18186528
$ws.Range("K126").Value = 42886956
$ws.Range("L126").Value = 54559584
$ws.Range("M126").Value = -42884486
$ws.Range("N126").Value = -54564524

# LTW row 59
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 149499.5
$ws.Range("J59").Value = 149499.5
$ws.Range("L59").Value = 149499.5
$ws.Range("N59").Value = -150807.5

# LTW row 100
$ws.Range("H100").Value = 2966.6
$ws.Range("I100").Value = 2303.3333
$ws.Range("J100").Value = 3509.2727
$ws.Range("K100").Value = 2303.3333
$ws.Range("L100").Value = 3509.2727
$ws.Range("M100").Value = -1762.3333
$ws.Range("N100").Value = -4591.2727

# LTW row 114
$ws.Range("H114").Value = 50198
$ws.Range("J114").Value = 50198
$ws.Range("L114").Value = 50198
$ws.Range("N114").Value = -58876

# LTW row 122
$ws.Range("H122").Value = 6678.4287
$ws.Range("I122").Value = 7187.25
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 21561.75
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -19111.75
$ws.Range("N122").Value = -22900
